# WTW Quarterly Financials update
# -----------------------------------------------------------------
# The source spreadsheet reports trailing quarters across columns.
# Two new quarters (Q ending 2018-12-29 and 2018-09-29) were added,
# which in Excel terms means inserting two new columns immediately
# before the existing "most recent quarter" column (D), pushing all
# prior quarters two columns to the right (old D:K -> new F:M) and
# then filling the new D:E columns with the latest figures across
# all three statements (Income Statement, Balance Sheet, Cash Flow
# Statement). One historical row (Capital Expenditures, row 91) was
# also restated with corrected figures for several of the shifted
# quarters.
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column D. Excel shifts the existing
# D:K data (and its formatting) to F:M automatically.
$ws.Range("D:E").Insert()

# Per-row new-quarter values for the two freshly inserted columns
# (D = quarter ending 2018-12-29, E = quarter ending 2018-09-29).
# IsDate marks the three "Period Ending" header rows (one per
# statement) which use the date number format instead of the
# thousands-separated number format used by every other data row.
$rows = @(
    @{ Row=7; D=43463; E=43372; IsDate=$true },
    @{ Row=8; D=330400; E=365800; IsDate=$false },
    @{ Row=9; D=145200; E=150400; IsDate=$false },
    @{ Row=10; D=185200; E=215400; IsDate=$false },
    @{ Row=11; D=$null; E=$null; IsDate=$false },
    @{ Row=12; D="NA"; E="NA"; IsDate=$false },
    @{ Row=13; D=0; E=0; IsDate=$false },
    @{ Row=14; D=0; E=0; IsDate=$false },
    @{ Row=15; D=0; E=0; IsDate=$false },
    @{ Row=16; D=$null; E=$null; IsDate=$false },
    @{ Row=17; D=250000; E=246900; IsDate=$false },
    @{ Row=18; D=80400; E=118900; IsDate=$false },
    @{ Row=19; D=$null; E=$null; IsDate=$false },
    @{ Row=20; D=-600; E=-900; IsDate=$false },
    @{ Row=21; D=91200; E=128800; IsDate=$false },
    @{ Row=22; D=35100; E=35500; IsDate=$false },
    @{ Row=23; D=44600; E=82500; IsDate=$false },
    @{ Row=24; D=3600; E=12400; IsDate=$false },
    @{ Row=25; D=0; E=0; IsDate=$false },
    @{ Row=26; D=41000; E=70100; IsDate=$false },
    @{ Row=27; D=41100; E=70100; IsDate=$false },
    @{ Row=28; D=0; E=0; IsDate=$false },
    @{ Row=29; D=2700; E="NA"; IsDate=$false },
    @{ Row=30; D=0; E=0; IsDate=$false },
    @{ Row=31; D=0; E=0; IsDate=$false },
    @{ Row=32; D=600; E=900; IsDate=$false },
    @{ Row=33; D=43800; E=70100; IsDate=$false },
    @{ Row=34; D=0; E=0; IsDate=$false },
    @{ Row=35; D=43800; E=70100; IsDate=$false },
    @{ Row=38; D=43463; E=43372; IsDate=$true },
    @{ Row=39; D=$null; E=$null; IsDate=$false },
    @{ Row=40; D=$null; E=$null; IsDate=$false },
    @{ Row=41; D=237000; E=219800; IsDate=$false },
    @{ Row=42; D=0; E=0; IsDate=$false },
    @{ Row=43; D=27200; E=28900; IsDate=$false },
    @{ Row=44; D=25900; E=14300; IsDate=$false },
    @{ Row=45; D=76400; E=73200; IsDate=$false },
    @{ Row=46; D=366400; E=336100; IsDate=$false },
    @{ Row=47; D=0; E=0; IsDate=$false },
    @{ Row=48; D=52200; E=49800; IsDate=$false },
    @{ Row=49; D=960800; E=962000; IsDate=$false },
    @{ Row=50; D=0; E=0; IsDate=$false },
    @{ Row=51; D=0; E=0; IsDate=$false },
    @{ Row=52; D=35100; E=33500; IsDate=$false },
    @{ Row=53; D=0; E=0; IsDate=$false },
    @{ Row=54; D=1414500; E=1381500; IsDate=$false },
    @{ Row=55; D=$null; E=$null; IsDate=$false },
    @{ Row=56; D=$null; E=$null; IsDate=$false },
    @{ Row=57; D=27100; E=23000; IsDate=$false },
    @{ Row=58; D=77000; E=57800; IsDate=$false },
    @{ Row=59; D=237200; E=231300; IsDate=$false },
    @{ Row=60; D=341300; E=312000; IsDate=$false },
    @{ Row=61; D=1669700; E=1687500; IsDate=$false },
    @{ Row=62; D=208500; E=223300; IsDate=$false },
    @{ Row=63; D=0; E=0; IsDate=$false },
    @{ Row=64; D=0; E=0; IsDate=$false },
    @{ Row=65; D=0; E=0; IsDate=$false },
    @{ Row=66; D=2223500; E=2226700; IsDate=$false },
    @{ Row=67; D=$null; E=$null; IsDate=$false },
    @{ Row=68; D=0; E=0; IsDate=$false },
    @{ Row=69; D=0; E=0; IsDate=$false },
    @{ Row=70; D=0; E=0; IsDate=$false },
    @{ Row=71; D=0; E=0; IsDate=$false },
    @{ Row=72; D=2382400; E=2340300; IsDate=$false },
    @{ Row=73; D=0; E=0; IsDate=$false },
    @{ Row=74; D=0; E=0; IsDate=$false },
    @{ Row=75; D=0; E=0; IsDate=$false },
    @{ Row=76; D=-808900; E=-845200; IsDate=$false },
    @{ Row=77; D=0; E=0; IsDate=$false },
    @{ Row=80; D=43463; E=43372; IsDate=$true },
    @{ Row=81; D=43800; E=70100; IsDate=$false },
    @{ Row=82; D=$null; E=$null; IsDate=$false },
    @{ Row=83; D=11500; E=10800; IsDate=$false },
    @{ Row=84; D=0; E=0; IsDate=$false },
    @{ Row=85; D=0; E=0; IsDate=$false },
    @{ Row=86; D=0; E=0; IsDate=$false },
    @{ Row=87; D=0; E=0; IsDate=$false },
    @{ Row=88; D=0; E=0; IsDate=$false },
    @{ Row=89; D=41000; E=86100; IsDate=$false },
    @{ Row=90; D=$null; E=$null; IsDate=$false },
    @{ Row=91; D=-7100; E=-4200; IsDate=$false },
    @{ Row=92; D=0; E=0; IsDate=$false },
    @{ Row=93; D=0; E=0; IsDate=$false },
    @{ Row=94; D=-19000; E=-14400; IsDate=$false },
    @{ Row=95; D=$null; E=$null; IsDate=$false },
    @{ Row=96; D=0; E=0; IsDate=$false },
    @{ Row=97; D=0; E=0; IsDate=$false },
    @{ Row=98; D=0; E=0; IsDate=$false },
    @{ Row=99; D=0; E=0; IsDate=$false },
    @{ Row=100; D=-3600; E=-19300; IsDate=$false },
    @{ Row=101; D=-1100; E=-300; IsDate=$false },
    @{ Row=102; D=17200; E=52000; IsDate=$false }

)

foreach ($item in $rows) {
    $r = $item.Row
    $dst = $ws.Range("D" + $r + ":E" + $r)

    if ($item.IsDate) {
        $dst.NumberFormat = "[$-409]d\-mmm\-yy;@"
    } else {
        $dst.NumberFormat = "#,##0"
        $dst.HorizontalAlignment = -4152   # xlRight
    }
    $dst.Font.Name = "Verdana"
    $dst.Font.Size = 12
    $dst.Font.Bold = $item.IsDate

    if ($item.D -ne $null) {
        $ws.Range("D" + $r).Value2 = $item.D
    }
    if ($item.E -ne $null) {
        $ws.Range("E" + $r).Value2 = $item.E
    }
}

# Row 91 ("Capital Expenditures") was also restated: besides the two
# new quarters (D91/E91, set above), the six quarters that shifted
# into F91:K91 (formerly D91:I91) carry corrected figures too. The
# last two shifted quarters (L91/M91, formerly J91/K91) are unchanged.
$row91 = @{ F = -6000; G = -1800; H = -3000; I = -5100; J = -2200; K = -10500 }
foreach ($col in $row91.Keys) {
    $ws.Range($col + "91").Value2 = $row91[$col]
}
